# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sat Apr  1 05:37:09 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.598.61"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").Value = "1.825.45"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'316.37"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "'0.5303"
$ws.Range("E7").Value = "  -2.64%  "

$ws.Range("D8").Value = "'0.3973"
$ws.Range("E8").Value = "  +4.55%  "

$ws.Range("D9").Value = "'0.07734"
$ws.Range("E9").Value = "  +3.34%  "

$ws.Range("D10").Value = "'42.08"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").Value = "'1.117"
$ws.Range("E11").Value = "  +1.90%  "

$ws.Range("D12").Value = "'21.10"
$ws.Range("E12").Value = "  +2.53%  "

$ws.Range("D13").Value = "'6.317"
$ws.Range("E13").Value = "  +1.59%  "

$ws.Range("D14").Value = "'1.003"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("D15").Value = "'7.568"
$ws.Range("E15").Value = "  +2.74%  "

$ws.Range("D16").Value = "1.826.83"
$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").Value = "'93.22"
$ws.Range("E17").Value = "  +3.55%  "

$ws.Range("E18").Value = "  +2.12%  "

$ws.Range("D19").Value = "'0.06617"
$ws.Range("E19").Value = "  +1.38%  "

$ws.Range("D20").Value = "'17.80"
$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "'6.084"
$ws.Range("E22").Value = "  +2.47%  "

$ws.Range("D23").Value = "28.606.92"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").Value = "'2.235"
$ws.Range("E25").Value = "  +6.81%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.063.24"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.78"
$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'156.79"
$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("D29").Value = "'2.412"
$ws.Range("E29").Value = "  +2.75%  "

$ws.Range("D30").Value = "'125.45"
$ws.Range("E30").Value = "  +2.53%  "

$ws.Range("D31").Value = "'1.150"
$ws.Range("E31").Value = "  +2.54%  "

$ws.Range("D32").Value = "'0.1125"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").Value = "'5.734"
$ws.Range("E33").Value = "  +2.64%  "

$ws.Range("D34").Value = "'3.662"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("D35").Value = "'0.07326"
$ws.Range("E35").Value = "  +4.72%  "

$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").Value = "'0.02351"
$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("D38").Value = "'8.907"
$ws.Range("E38").Value = "  +4.86%  "

$ws.Range("D39").Value = "'5.197"
$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("D40").Value = "'11.40"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").Value = "'0.6291"
$ws.Range("E41").Value = "  +1.55%  "

$ws.Range("D42").Value = "'1.196"
$ws.Range("E42").Value = "  +1.86%  "

$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "'1.398"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("D45").Value = "'13.57"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("D46").Value = "'0.5935"
$ws.Range("E46").Value = "  +2.94%  "

$ws.Range("E47").Value = "  +1.04%  "

$ws.Range("D48").Value = "'125.53"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").Value = "'2.000"
$ws.Range("E49").Value = "  +3.80%  "

$ws.Range("D50").Value = "'1.192"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").Value = "'0.06955"
$ws.Range("E51").Value = "  +1.86%  "
